# Generate Report for Handback
#
# 9438add3-dae0-4160-81c9-9872e95225a0.md and d7691be1-dd41-4ca5-86fa-28f176535f88.md
# finished their zh-cn / de-de handback cycle. Flip their status from
# "Ready for handoff" to "Handed back: in sync with en-US" on the Overview
# sheet, and on the per-locale sheets populate the "Latest Target File" /
# "Latest Handback File" hyperlinks plus the new "Latest Handback DateTime".

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: rows 3 (9438add3...) and 4 (d7691be1...) move from
# "Ready for handoff" to "Handed back: in sync with en-US" in both the
# zh-cn (B) and de-de (C) columns. Column D (Latest Handoff Date) is
# unchanged.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

$wsOverview.Range("B4").Value = $statusHandedBack
$wsOverview.Range("C4").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet: rows 3 & 4 (both describe 9438add3...zh-cn.xlf) get their
# Status flipped, the Latest Target File / Latest Handback File
# hyperlinks filled in, and the Latest Handback DateTime stamped.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("C4").Value = $statusHandedBack

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/739575f8a230f98065512eed7a4977413ed8b303/e2e/9438add3-dae0-4160-81c9-9872e95225a0.md", "", "", "9438add3-dae0-4160-81c9-9872e95225a0.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/739575f8a230f98065512eed7a4977413ed8b303/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/9438add3-dae0-4160-81c9-9872e95225a0.d64e3cf96924ff0577bd5b31669f0de2eb2a9295.zh-cn.xlf", "", "", "9438add3-dae0-4160-81c9-9872e95225a0.d64e3cf96924ff0577bd5b31669f0de2eb2a9295.zh-cn.xlf") | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/739575f8a230f98065512eed7a4977413ed8b303/e2e/9438add3-dae0-4160-81c9-9872e95225a0.md", "", "", "9438add3-dae0-4160-81c9-9872e95225a0.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/739575f8a230f98065512eed7a4977413ed8b303/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/9438add3-dae0-4160-81c9-9872e95225a0.d64e3cf96924ff0577bd5b31669f0de2eb2a9295.zh-cn.xlf", "", "", "9438add3-dae0-4160-81c9-9872e95225a0.d64e3cf96924ff0577bd5b31669f0de2eb2a9295.zh-cn.xlf") | Out-Null

$wsZhCn.Range("H3").Value = "2016-03-19 03:42:08"
$wsZhCn.Range("H4").Value = "2016-03-19 03:42:08"

# ---------------------------------------------------------------------
# de-de sheet: same shape of change as zh-cn, but the datetime stamp is
# different (2016-03-19 03:42:22) and file names carry the de-de suffix.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("C4").Value = $statusHandedBack

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/e680230e96df8fbb7bfaa52f4f1ff42ec0bd43be/e2e/9438add3-dae0-4160-81c9-9872e95225a0.md", "", "", "9438add3-dae0-4160-81c9-9872e95225a0.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e680230e96df8fbb7bfaa52f4f1ff42ec0bd43be/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/9438add3-dae0-4160-81c9-9872e95225a0.d64e3cf96924ff0577bd5b31669f0de2eb2a9295.de-de.xlf", "", "", "9438add3-dae0-4160-81c9-9872e95225a0.d64e3cf96924ff0577bd5b31669f0de2eb2a9295.de-de.xlf") | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/e680230e96df8fbb7bfaa52f4f1ff42ec0bd43be/e2e/9438add3-dae0-4160-81c9-9872e95225a0.md", "", "", "9438add3-dae0-4160-81c9-9872e95225a0.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e680230e96df8fbb7bfaa52f4f1ff42ec0bd43be/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/9438add3-dae0-4160-81c9-9872e95225a0.d64e3cf96924ff0577bd5b31669f0de2eb2a9295.de-de.xlf", "", "", "9438add3-dae0-4160-81c9-9872e95225a0.d64e3cf96924ff0577bd5b31669f0de2eb2a9295.de-de.xlf") | Out-Null

$wsDeDe.Range("H3").Value = "2016-03-19 03:42:22"
$wsDeDe.Range("H4").Value = "2016-03-19 03:42:22"

Write-Output "Handback report generated"
